$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts the existing rows 10-48
# (the rest of the "Femacal de La Calera" Arveja Verde series) down to 11-49.
$ws.Rows.Item(10).EntireRow.Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Femacal de La Calera"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44560
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 100112022
$ws.Range("G10").Value = "Arveja Verde"
$ws.Range("H10").Value = "Perfection"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 27000
$ws.Range("L10").Value = 28000
$ws.Range("M10").Value = 27500
$ws.Range("N10").Value = "$/malla 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 1100
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
